$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange
$lastRealParagraph = $tr.Paragraphs(2, 1)
[void]$lastRealParagraph.InsertAfter([char]13 + "https://github.com/robojay/SOAR")
